# Applies the scheduled market-data refresh to the per-job Leve profit sheets.
# For each changed row, currentAveragePrice/NQ/HQ, LevePriceNQ/HQ and the
# resulting LeveProfitNQ/HQ columns (H:N) are updated to the newly-fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2001.1428
$ws.Range("I62").Value = 1702
$ws.Range("K62").Value = 1702
$ws.Range("M62").Value = -1078

$ws.Range("H65").Value = 2001.1428
$ws.Range("I65").Value = 1702
$ws.Range("K65").Value = 8510
$ws.Range("M65").Value = -5390

$ws.Range("H82").Value = 4042
$ws.Range("I82").Value = 4042
$ws.Range("K82").Value = 12126
$ws.Range("M82").Value = -11720

$ws.Range("H85").Value = 4042
$ws.Range("I85").Value = 4042
$ws.Range("K85").Value = 12126
$ws.Range("M85").Value = -10722

$ws.Range("H98").Value = 1554.3125
$ws.Range("I98").Value = 348
$ws.Range("K98").Value = 348
$ws.Range("M98").Value = 1150

$ws.Range("H100").Value = 2597.5833
$ws.Range("I100").Value = 2117.3
$ws.Range("J100").Value = 4999
$ws.Range("K100").Value = 2117.3
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -1576.3
$ws.Range("N100").Value = -6081

$ws.Range("H122").Value = 1554.3125
$ws.Range("I122").Value = 348
$ws.Range("K122").Value = 1044
$ws.Range("M122").Value = 1406

$ws.Range("H132").Value = 14387.777
$ws.Range("I132").Value = 14387.777
$ws.Range("K132").Value = 43163.331
$ws.Range("M132").Value = -40633.331

$ws.Range("H135").Value = 35715530
$ws.Range("I135").Value = 1439.8
$ws.Range("J135").Value = 125000750
$ws.Range("K135").Value = 12958.2
$ws.Range("L135").Value = 1125006750
$ws.Range("M135").Value = -10423.2
$ws.Range("N135").Value = -1125011820

$ws.Range("H138").Value = 5227.1665
$ws.Range("I138").Value = 2629
$ws.Range("J138").Value = 5463.364
$ws.Range("K138").Value = 7887
$ws.Range("L138").Value = 16390.092
$ws.Range("M138").Value = -2747
$ws.Range("N138").Value = -26670.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3537.25
$ws.Range("I2").Value = 1414.8572
$ws.Range("K2").Value = 1414.8572
$ws.Range("M2").Value = -1301.8572

$ws.Range("H32").Value = 4415.58
$ws.Range("I32").Value = 3928.5112
$ws.Range("K32").Value = 3928.5112
$ws.Range("M32").Value = -3641.5112

$ws.Range("H45").Value = 7558.4287
$ws.Range("I45").Value = 7482
$ws.Range("K45").Value = 7482
$ws.Range("M45").Value = -7105

$ws.Range("H74").Value = 58895372
$ws.Range("I74").Value = 77015680
$ws.Range("J74").Value = 4375
$ws.Range("K74").Value = 77015680
$ws.Range("L74").Value = 4375
$ws.Range("M74").Value = -77014806
$ws.Range("N74").Value = -6123

$ws.Range("H77").Value = 58895372
$ws.Range("I77").Value = 77015680
$ws.Range("J77").Value = 4375
$ws.Range("K77").Value = 385078400
$ws.Range("L77").Value = 21875
$ws.Range("M77").Value = -385074032
$ws.Range("N77").Value = -30611

$ws.Range("H116").Value = 3537.25
$ws.Range("I116").Value = 1414.8572
$ws.Range("K116").Value = 1414.8572
$ws.Range("M116").Value = 879.1428000000001

$ws.Range("H122").Value = 7753667.5
$ws.Range("I122").Value = 1666.5358
$ws.Range("J122").Value = 22224070
$ws.Range("K122").Value = 4999.607400000001
$ws.Range("L122").Value = 66672210
$ws.Range("M122").Value = -2549.607400000001
$ws.Range("N122").Value = -66677110

$ws.Range("H132").Value = 32314222
$ws.Range("I132").Value = 13412.958
$ws.Range("K132").Value = 40238.874
$ws.Range("M132").Value = -37708.874

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3537.25
$ws.Range("I3").Value = 1414.8572
$ws.Range("K3").Value = 1414.8572
$ws.Range("M3").Value = -1300.8572

$ws.Range("H26").Value = 19418.084
$ws.Range("I26").Value = 18001.545
$ws.Range("K26").Value = 18001.545
$ws.Range("M26").Value = -17709.545

$ws.Range("H94").Value = 1468.1428
$ws.Range("I94").Value = 1035.4166
$ws.Range("J94").Value = 2045.1111
$ws.Range("K94").Value = 1035.4166
$ws.Range("L94").Value = 2045.1111
$ws.Range("M94").Value = -584.4166
$ws.Range("N94").Value = -2947.1111

$ws.Range("H96").Value = 28913.732
$ws.Range("I96").Value = 6250.636
$ws.Range("J96").Value = 91237.25
$ws.Range("K96").Value = 6250.636
$ws.Range("L96").Value = 91237.25
$ws.Range("M96").Value = -3504.636
$ws.Range("N96").Value = -96729.25

$ws.Range("H134").Value = 4379.2173
$ws.Range("I134").Value = 4148.3335
$ws.Range("K134").Value = 12445.0005
$ws.Range("M134").Value = -9910.000499999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7815688.5
$ws.Range("I31").Value = 2200.3333
$ws.Range("K31").Value = 2200.3333
$ws.Range("M31").Value = -1905.3333

$ws.Range("H34").Value = 7815688.5
$ws.Range("I34").Value = 2200.3333
$ws.Range("K34").Value = 2200.3333
$ws.Range("M34").Value = -1998.3333

$ws.Range("H58").Value = 3502.6
$ws.Range("I58").Value = 2500
$ws.Range("K58").Value = 2500
$ws.Range("M58").Value = -2297

$ws.Range("H132").Value = 94261.17999999999
$ws.Range("I132").Value = 94261.17999999999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 282783.54
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -280253.54
$ws.Range("N132").ClearContents() | Out-Null

$ws.Range("H134").Value = 3449.6
$ws.Range("I134").Value = 2374.6
$ws.Range("J134").Value = 5599.6
$ws.Range("K134").Value = 7123.799999999999
$ws.Range("L134").Value = 16798.8
$ws.Range("M134").Value = -4588.799999999999
$ws.Range("N134").Value = -21868.8

$ws.Range("H136").Value = 3502.6
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2546.9
$ws.Range("I5").Value = 310.8
$ws.Range("J5").Value = 4783
$ws.Range("K5").Value = 932.4000000000001
$ws.Range("L5").Value = 14349
$ws.Range("M5").Value = -820.4000000000001
$ws.Range("N5").Value = -14573

$ws.Range("H63").Value = 749.5
$ws.Range("J63").Value = 1000
$ws.Range("L63").Value = 3000
$ws.Range("N63").Value = -4498

$ws.Range("H66").Value = 749.5
$ws.Range("J66").Value = 1000
$ws.Range("L66").Value = 9000
$ws.Range("N66").Value = -16488

$ws.Range("H75").Value = 891.875
$ws.Range("I75").Value = 1078.6666
$ws.Range("J75").Value = 779.8
$ws.Range("K75").Value = 3235.9998
$ws.Range("L75").Value = 2339.4
$ws.Range("M75").Value = -2237.9998
$ws.Range("N75").Value = -4335.4

$ws.Range("H78").Value = 891.875
$ws.Range("I78").Value = 1078.6666
$ws.Range("J78").Value = 779.8
$ws.Range("K78").Value = 9707.999400000001
$ws.Range("L78").Value = 7018.2
$ws.Range("M78").Value = -4715.999400000001
$ws.Range("N78").Value = -17002.2

$ws.Range("H135").Value = 2546.9
$ws.Range("I135").Value = 310.8
$ws.Range("J135").Value = 4783
$ws.Range("K135").Value = 2797.2
$ws.Range("L135").Value = 43047
$ws.Range("M135").Value = -262.2000000000003
$ws.Range("N135").Value = -48117

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 229091.44
$ws.Range("I70").Value = 404432.4
$ws.Range("K70").Value = 404432.4
$ws.Range("M70").Value = -404162.4

$ws.Range("H73").Value = 229091.44
$ws.Range("I73").Value = 404432.4
$ws.Range("K73").Value = 404432.4
$ws.Range("M73").Value = -403496.4

$ws.Range("H102").Value = 6211
$ws.Range("I102").Value = 2091.2222
$ws.Range("K102").Value = 2091.2222
$ws.Range("M102").Value = -469.2222000000002

$ws.Range("H122").Value = 22728526
$ws.Range("I122").Value = 1187
$ws.Range("K122").Value = 3561
$ws.Range("M122").Value = -1111

$ws.Range("H132").Value = 2899
$ws.Range("J132").Value = 2499.6667
$ws.Range("L132").Value = 7499.000100000001
$ws.Range("N132").Value = -12559.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4162.8237
$ws.Range("I7").Value = 4199.6665
$ws.Range("J7").Value = 4121.375
$ws.Range("K7").Value = 4199.6665
$ws.Range("L7").Value = 4121.375
$ws.Range("M7").Value = -4087.6665
$ws.Range("N7").Value = -4345.375

$ws.Range("H122").Value = 2608597.8
$ws.Range("I122").Value = 3830.818
$ws.Range("K122").Value = 11492.454
$ws.Range("M122").Value = -9042.454000000002

$ws.Range("H126").Value = 4162.8237
$ws.Range("I126").Value = 4199.6665
$ws.Range("J126").Value = 4121.375
$ws.Range("K126").Value = 12598.9995
$ws.Range("L126").Value = 12364.125
$ws.Range("M126").Value = -10128.9995
$ws.Range("N126").Value = -17304.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7181643.5
$ws.Range("I122").Value = 64122.188
$ws.Range("K122").Value = 192366.564
$ws.Range("M122").Value = -189916.564

$ws.Range("H132").Value = 4091.5715
$ws.Range("I132").Value = 3210.375
$ws.Range("J132").Value = 5266.5
$ws.Range("K132").Value = 9631.125
$ws.Range("L132").Value = 15799.5
$ws.Range("M132").Value = -7101.125
$ws.Range("N132").Value = -20859.5

$ws.Range("H136").Value = 10707.962
$ws.Range("I136").Value = 4265.846
$ws.Range("K136").Value = 12797.538
$ws.Range("M136").Value = -10247.538
